$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.006.43"
$ws.Range("E2").Value = "  +8.61%  "

$ws.Range("D3").Value = "3.140.60"
$ws.Range("E3").Value = "  +5.99%  "

$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "589.50"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.88%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.76"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +8.25%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").Value = "3.138.00"
$ws.Range("E8").Value = "  +5.98%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.535"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.70%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.154"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +17.27%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.83"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +10.36%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.472"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +5.56%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000251"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +10.59%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.00"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +7.80%  "

$ws.Range("E15").Value = "  +1.28%  "

$ws.Range("D16").Value = "3.656.97"
$ws.Range("E16").Value = "  +5.88%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.20"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.43%  "

$ws.Range("D18").Value = "63.835.97"
$ws.Range("E18").Value = "  +8.10%  "

$ws.Range("D19").Value = "3.128.81"
$ws.Range("E19").Value = "  +5.54%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "474.09"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +9.51%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.25"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.63%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.734"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.41%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.60"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +9.15%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.37"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.79%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "82.52"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.73%  "

$ws.Range("E26").Value = "  +0.20%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.73"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +14.20%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.24"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.62%  "

$ws.Range("E29").Value = "  +6.34%  "

$ws.Range("E30").Value = "  -0.12%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.94"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +12.67%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "27.26"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.72%  "

$ws.Range("E33").Value = "  +6.93%  "

$ws.Range("D34").Value = "0.0₃0888"
$ws.Range("E34").Value = "  +17.32%  "

$ws.Range("E35").Value = "  +19.75%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.06"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +8.36%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.41"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +24.66%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.15"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.22%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "50.92"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.38%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "446.69"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +13.09%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.75"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.62%  "

$ws.Range("D42").Value = "2.941.64"
$ws.Range("E42").Value = "  +8.28%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0373"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +7.12%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.283"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +13.95%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.111"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +6.56%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.21"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +12.41%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "35.80"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +5.75%  "

$ws.Range("E48").Value = "  +0.01%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "123.72"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.32%  "

$ws.Range("E50").Value = "  +2.25%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "24.91"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +8.35%  "
